$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I9").Value = "If no new operation is required for 10s calculator shall be OFF"
